$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A26").Value = "JD_025"
$ws.Range("B26").Value = "Senior Associate Engineer"
$ws.Range("C26").Value = "Testing"
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = 4
